# Re-shuffle the per-occurrence data (columns A:AY) among rows 29-41 of the
# "Artfynd" sheet so that each row ends up holding the occurrence record
# that the commit moved it to. Rows 32, 35 and 40 are unaffected (they map
# to themselves). This mirrors the upstream source-data re-sync; the set of
# records is unchanged, only which row each one lives on.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# target row -> source row (i.e. "what ends up at row X used to be at row Y")
$mapping = @{
    29 = 34
    30 = 41
    31 = 33
    32 = 32
    33 = 37
    34 = 29
    35 = 35
    36 = 30
    37 = 31
    38 = 39
    39 = 36
    40 = 40
    41 = 38
}

$firstCol = 1          # A
$lastCol  = 51          # AY

# Columns that are genuinely numeric / boolean in this sheet; everything
# else in the occurrence block is plain text (even when it looks like a
# number, a date or a time - e.g. "1", "2023-09-17", "10:36").
$numericCols = @(1, 2, 5, 17, 18, 19)        # A, B, E, Q, R, S
$boolCols    = @(30, 31, 33)                 # AD, AE, AG

function Get-ColType($col) {
    if ($numericCols -contains $col) { return "n" }
    if ($boolCols -contains $col)    { return "b" }
    return "s"
}

# ---- 1. Snapshot every source row (values only, before anything is
#         overwritten) so the permutation can be applied safely even though
#         rows read from and rows written to overlap. ----
$snapshot = @{}
foreach ($targetRow in $mapping.Keys) {
    $srcRow = $mapping[$targetRow]
    if (-not $snapshot.ContainsKey($srcRow)) {
        $rowVals = @{}
        for ($col = $firstCol; $col -le $lastCol; $col++) {
            $cell = $ws.Cells.Item($srcRow, $col)
            $t = Get-ColType $col
            if ($t -eq "s") {
                $rowVals[$col] = [string]$cell.Text
                if ($cell.Text -eq "" -and $cell.Value2 -eq $null) {
                    $rowVals[$col] = $null
                }
            } else {
                $rowVals[$col] = $cell.Value2
            }
        }
        $snapshot[$srcRow] = $rowVals
    }
}

# ---- 2. Write every target row from the snapshot. ----
foreach ($targetRow in ($mapping.Keys | Sort-Object)) {
    $srcRow = $mapping[$targetRow]
    $rowVals = $snapshot[$srcRow]

    # Clear the whole occurrence block first so any column that has no
    # value in the source row ends up blank in the target row too.
    $ws.Range($ws.Cells.Item($targetRow, $firstCol), $ws.Cells.Item($targetRow, $lastCol)).ClearContents()

    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $val = $rowVals[$col]
        if ($val -eq $null) { continue }

        $cell = $ws.Cells.Item($targetRow, $col)
        $t = Get-ColType $col
        if ($t -eq "n") {
            $cell.Value2 = [double]$val
        } elseif ($t -eq "b") {
            if ($val -is [string]) {
                $cell.Value2 = ($val -eq "True" -or $val -eq "-1" -or $val -eq "1")
            } else {
                $cell.Value2 = [bool]$val
            }
        } else {
            if ($val -eq "") {
                $cell.Value = "'"
            } else {
                # Leading apostrophe forces text storage so numeric-looking
                # ("1"), date-looking ("2023-09-17") and time-looking
                # ("10:36") strings round-trip as text, not as a number /
                # serial date.
                $cell.Value = "'" + $val
            }
        }
    }
}

"Row data re-synced for rows 29-41."
